# Update the Price (D) and Volume(1h) (E) columns with refreshed crypto
# quote data. Values are leading-apostrophe-prefixed so Excel stores them
# as literal text (matching the workbook's existing text-based cells)
# instead of auto-converting numeric-looking / percent-looking strings
# into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''298.36'
$ws.Range("E2").Value = '''-2.34%'
$ws.Range("D3").Value = '''31.27'
$ws.Range("E3").Value = '''-2.64%'
$ws.Range("D4").Value = '''5.102'
$ws.Range("E4").Value = '''-2.33%'
$ws.Range("D5").Value = '''0.07898'
$ws.Range("E5").Value = '''5.60%'
$ws.Range("D6").Value = '''2.336'
$ws.Range("E6").Value = '''-1.12%'
$ws.Range("D7").Value = '''7.742'
$ws.Range("E7").Value = '''-3.18%'
$ws.Range("D8").Value = '''3.863'
$ws.Range("E8").Value = '''-0.29%'
$ws.Range("D9").Value = '''0.9214'
$ws.Range("E9").Value = '''0.42%'
$ws.Range("D10").Value = '''0.1726'
$ws.Range("E10").Value = '''-0.74%'
$ws.Range("D11").Value = '''0.07355'
$ws.Range("E11").Value = '''-4.80%'
$ws.Range("D12").Value = '''0.08962'
$ws.Range("E12").Value = '''8.69%'
$ws.Range("D13").Value = '''0.03008'
$ws.Range("E13").Value = '''0.03%'
$ws.Range("D14").Value = '''0.1001'
$ws.Range("E14").Value = '''0.44%'
$ws.Range("D15").Value = '''0.001509'
$ws.Range("E15").Value = '''0.40%'
$ws.Range("D16").Value = '''0.006145'
$ws.Range("E16").Value = '''-0.15%'
$ws.Range("D17").Value = '''3.486'
$ws.Range("E17").Value = '''-0.19%'
$ws.Range("E18").Value = '''2.78%'
$ws.Range("E19").Value = '''0.20%'
$ws.Range("E20").Value = '''-1.42%'
$ws.Range("D21").Value = '''4.160'
$ws.Range("E21").Value = '''-10.47%'
$ws.Range("E22").Value = '''8.89%'
$ws.Range("D23").Value = '''0.04619'
$ws.Range("E23").Value = '''0.36%'
$ws.Range("D24").Value = '''0.001248'
$ws.Range("E24").Value = '''0.25%'
$ws.Range("D25").Value = '''0.004460'
$ws.Range("E25").Value = '''-1.54%'
$ws.Range("D26").Value = '''0.0001199'
$ws.Range("E26").Value = '''-7.42%'
$ws.Range("D27").Value = '''0.0003394'
$ws.Range("E27").Value = '''24.15%'
$ws.Range("D39").Value = '''0.01740'
$ws.Range("E39").Value = '''-2.65%'
$ws.Range("D40").Value = '''0.04590'
$ws.Range("E40").Value = '''0.27%'
$ws.Range("D41").Value = '''0.006970'
$ws.Range("E41").Value = '''-5.48%'
$ws.Range("E42").Value = '''-0.56%'
$ws.Range("D43").Value = '''0.002188'
$ws.Range("E43").Value = '''0.76%'
$ws.Range("D44").Value = '''0.009562'
$ws.Range("E44").Value = '''-11.46%'
$ws.Range("D45").Value = '''0.00006269'
$ws.Range("E45").Value = '''-3.49%'
$ws.Range("D46").Value = '''0.00000000749'
$ws.Range("E46").Value = '''-0.10%'
$ws.Range("D47").Value = '''0.007972'
$ws.Range("E47").Value = '''-19.23%'
$ws.Range("D48").Value = '''0.7476'
$ws.Range("E48").Value = '''-8.89%'
$ws.Range("D49").Value = '''0.00002098'
$ws.Range("E49").Value = '''-0.10%'
$ws.Range("D50").Value = '''0.0001998'
$ws.Range("E50").Value = '''-0.03%'
